$d = $word.ActiveDocument

# --- Paragraph 1: date change 30.07.24 -> 29.07.24 ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)
$r1.Text = '⚡️🚀המאמר היומי של מייק 29.07.24: ⚡️🚀'

# --- Paragraph 2: title change, remove the manual line break run ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
[void]$r2.MoveEnd(1, -1)
$r2.Text = 'Large Scale Dataset Distillation with Domain Shift'

# --- Paragraph 3 ---
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
[void]$r3.MoveEnd(1, -1)
$r3.Text = 'המאמר מציע שיטה מעניינת ודי מקורית לגנרוט דאטה מהתפלגות הנתונה על ידי דאטהסט מתויג. למשל בהינתן דאטהסט של תמונות D_s המטרה היא ליצור דאטהסט (מתויג) גדול בעל התפלגות ה"מושרה" על ידי D_s. המחברים טוענים כי השיטות הקיימות מתקשות לבנות(distill) דאטהסט גדול המשקף בצורה נאמנה את המאפיינים המהותיים של D_s.'

# --- Paragraph 4 ---
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
[void]$r4.MoveEnd(1, -1)
$r4.Text = 'המחברים מציעים לגשת לבעיה זו עם גישה מעולם של domain adaption או DA בקצרה. בגדול מאוד DA  היא תהליך של "התאמת מודל" במקרים בהם התפלגות הדאטה בזמן האינפרנס שונה מזו של הדאטה שעליה אומן המודל. התחום הזה עשיר בשיטות שחלקן די מורכבות מתמטיות ומערבות לרוב מינימיזציה של מרחק בין התפלגויות הדאטה (KL וכאלה).'

# --- Paragraph 5 ---
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
[void]$r5.MoveEnd(1, -1)
$r5.Text = 'למעשה המאמר המסוקר מתרגם את בעיית יצירת הדאטה לבעיית DA. התפלגות הדאטהסט שאנו מגנרטים ״ממנו״ D_s משחק תפקיד של התפלגות המקור במקרה של DA (שעליו מאומן המודל ב-DA) ואילו התפלגות הדאטה המגונרט משחקת תפקיד של התפלגות היעד D_t (כלומר זו של הדאטה שעליו מפעילים את המודל ב-DA). המטרה כאן לאמן מודל המקרב את ההתפלגויות האל.'

# --- Paragraph 6 ---
$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
[void]$r6.MoveEnd(1, -1)
$r6.Text = 'אבל איך נחשב את ההתפלגויות האלו? המאמר מייצג את ההתפלגויות האלו על ידי התפלגות של האקטיבציות של השכבות השונות של הרשת. בפשטות עבור הדאטסט D_s אנו מייצגים את התפלגות הדאטה על ידי וקטור הממוצעים ומטריצת קווריאנס של כל השכבות של המודל M_s(מניחים שהם גאוסיים). בדיוק באותו האופן אנו מייצגים את ההתפלגות של הדאטה המגונרט. '

# --- Paragraph 7 ---
$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
[void]$r7.MoveEnd(1, -1)
$r7.Text = 'אבל מה כאן M_s ומה עושים כדי לקרב את התפלגות של הדאטה המגונרט להתפלגות הדאטה האמיתי? המודל M_s אומן לשערך את ההתפלגות של הדאטהסט המתויג D_s (המאמר לא מפרט איך M_s מאומן בדיוק). למעשה האופטימיזציה מתבצעת על הדאטה המגונרט כאשר המודל M_s נותר ללא שינוי. כלומר מתחילים מתמונות הנדגמות באקראי עם הלייבלים והמטרה היא לבצע מורד הגרדיאנט(gradient descent) על התמונות האלו במטרה לקרב אותם להתפלגות של D_s.'

# --- Paragraph 8 ---
$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
[void]$r8.MoveEnd(1, -1)
$r8.Text = 'עכשיו נשאלת השאלה מפונקציית הלוס כאן. כאמור בשלב הראשון אנו מאפטמים את התמונות המגונרטות במטרה למזער מרחק KL בין התפלגויות המשקלי המודל M_s(נותר ללא שינוי) של D_s (נותר קבוע לכל אורך הדרך) ולבין התפלגות של משקלי המודל M_s עבור D_t. המחברים מניחים ששתי התפלגויות אלו הם גאוסיים שעבורם מרחק KL ניתן לחישוב באופן מדויק בהינתם וקטורי תוחלות ומטריצות קווריאנס של D_s ו-  D_t עם M_s. איבר נוסף בלוס מנסה למקסם (=למזער עם סימן מינוס) הוא ההתפלגות המותנית של לייבל y בהינתן פיסת דאטה מג''ונרט (הרי אנו מגנרטים דאטה מתיוג). התיוג של כל פיסת דאטה מגונרטת נקבע מראש ולא משתנה במהלך האימון.'

# --- Append 3 new paragraphs (9, 10, 11) after paragraph 8 ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
[void]$lastRange.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = 'השלב השני הוא מזעור של מרחק KL בין ההתפלגות המותנית של הלייבלים של הדאטה המגונרט לבין זה של הדאטה מ-D_s. בשביל כך מנצלים את הדאטה המגונרט מהשלב הראשון. מחשבים את התפלגות הלייבלים עבור הדאטה המגונרט הזה עם מודל M_s ומאפטמים את הדאטה המגונרט במטרה לקרב את שתי ההתפלגויות האלו של הלייבלים (של הדאטה המגונרט ושל הדאטה מ-D_s).'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
[void]$lastRange.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = 'יש עוד לא מעט פרטים מעניינים על איך בדיוק מתבצע האימון (משתמשים בלא מעט מודלים לחישוב סטטיסטיקות המשקלים, עושים מיצוע מעריכים לסטטיסטיקות של הבאצ''ים וכדומה). המאמר לא כתוב מאוד ברור אבל הרעיון יפה.'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
[void]$lastRange.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = 'https://dl.acm.org/doi/10.5555/3692070.3693400'

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
